$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

# Copy formatting (including the date/time number-format style) from the
# row above so the new date cell reuses the existing style instead of a
# brand new one being created.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 45450.2916666667

$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 2
$ws.Cells.Item($row, 4).Value = 2
$ws.Cells.Item($row, 5).Value = 2
$ws.Cells.Item($row, 6).Value = 2

# Store "2" as text (shared string), matching the source data, rather
# than letting Excel coerce it into a number.
$ws.Cells.Item($row, 7).NumberFormat = "@"
$ws.Cells.Item($row, 7).Value = "2"
$ws.Cells.Item($row, 7).Style = "Normal"

$ws.Cells.Item($row, 8).Value = "KK.MI"
